$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new first column (shifts B->C, C->D, D->E, A->B) -------------
# This preserves the existing column widths/styles of the shifted columns
# exactly, matching what Excel itself does on Insert().
$ws.Columns.Item(1).Insert()

# --- New "TabName" / "CasesTab" column -------------------------------------
$ws.Range("A1").Value = "TabName"
$ws.Range("A2").Value = "CasesTab"
# Target best-fit width is ~8.8164 chars; ColumnWidth is requested slightly
# narrower so the host's own char/pixel rounding lands on the closest
# representable stored width.
$ws.Columns.Item(1).ColumnWidth = 7.983072916666667

# --- Updated WebExcel query (now aliases clinical_trial as `ct`) -----------
$webQuery = @'
MATCH (ct:clinical_trial)<--(a:arm)<--(c:case)
    WHERE c.ethnicity IN ['HISPANIC_OR_LATINO']
WITH DISTINCT c, a, ct
RETURN 
    COALESCE(c.case_id, '') AS `Case ID`,
    COALESCE(ct.clinical_trial_designation, '') AS `Trial Code`,
    COALESCE(a.arm_id, '') AS `Arm`,
    COALESCE(a.arm_drug, '') AS `Arm Treatment`,
    COALESCE(c.disease, '') AS `Diagnosis`,
    COALESCE(c.gender, '') AS `Gender`,
    COALESCE(c.race, '') AS `Race`,
    COALESCE(c.ethnicity, '') AS `Ethnicity`
'@
$ws.Range("B2").Value = $webQuery

# --- Updated StatQuery query (specimen-rooted traversal) --------------------
$statQuery = @'
MATCH (s:specimen)-->(c:case)-->(:arm)-->(ct:clinical_trial)
    WHERE WHERE c.ethnicity IN ['HISPANIC_OR_LATINO']
OPTIONAL MATCH (f:file)-->(:sequencing_assay)-->(:nucleic_acid)-->(s)
RETURN 
    COUNT(DISTINCT f) AS number_of_files,
    COUNT(DISTINCT c.case_id) AS number_of_cases,
    COUNT(DISTINCT ct.clinical_trial_designation) AS number_of_trials
'@
$ws.Range("C2").Value = $statQuery

# --- Row height for the now-taller wrapped query text -----------------------
$ws.Rows.Item(2).RowHeight = 174

Write-Host "done"
